# Updated cryptos list on Tue Nov  7 17:09:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.666.16"
$ws.Range("E2").Value = "'  -1.88%  "
$ws.Range("D3").Value = "'1.855.68"
$ws.Range("E3").Value = "'  -2.90%  "
$ws.Range("E4").Value = "'  -0.83%  "
$ws.Range("D5").Value = "'242.64"
$ws.Range("E5").Value = "'  -4.43%  "
$ws.Range("D6").Value = "'0.663"
$ws.Range("E6").Value = "'  -8.21%  "
$ws.Range("E7").Value = "'  -0.91%  "
$ws.Range("D8").Value = "'41.35"
$ws.Range("E8").Value = "'  +1.88%  "
$ws.Range("D9").Value = "'0.336"
$ws.Range("E9").Value = "'  -6.18%  "
$ws.Range("D10").Value = "'0.0721"
$ws.Range("E10").Value = "'  -4.13%  "
$ws.Range("D11").Value = "'0.0964"
$ws.Range("E11").Value = "'  -2.81%  "
$ws.Range("B12").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "'2.126.06"
$ws.Range("E12").Value = "'  -2.95%  "
$ws.Range("B13").Value = "'Chainlink"
$ws.Range("C13").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'12.63"
$ws.Range("E13").Value = "'  +0.16%  "
$ws.Range("D14").Value = "'0.698"
$ws.Range("E14").Value = "'  -2.58%  "
$ws.Range("D15").Value = "'1.856.03"
$ws.Range("E15").Value = "'  -3.00%  "
$ws.Range("D16").Value = "'4.74"
$ws.Range("E16").Value = "'  -3.63%  "
$ws.Range("D17").Value = "'34.643.00"
$ws.Range("E17").Value = "'  -1.97%  "
$ws.Range("D18").Value = "'71.38"
$ws.Range("E18").Value = "'  -4.13%  "
$ws.Range("D19").Value = "'0.0₃0798"
$ws.Range("E19").Value = "'  -6.05%  "
$ws.Range("D20").Value = "'240.33"
$ws.Range("E20").Value = "'  -1.41%  "
$ws.Range("D21").Value = "'12.36"
$ws.Range("E21").Value = "'  -5.02%  "
$ws.Range("D22").Value = "'4.80"
$ws.Range("E22").Value = "'  -5.42%  "
$ws.Range("E23").Value = "'  -0.98%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("E24").Value = "'  +4.07%  "
$ws.Range("D25").Value = "'2.11"
$ws.Range("E25").Value = "'  -14.39%  "
$ws.Range("D26").Value = "'162.48"
$ws.Range("E26").Value = "'  -2.52%  "
$ws.Range("D27").Value = "'8.20"
$ws.Range("E27").Value = "'  -5.05%  "
$ws.Range("D28").Value = "'17.83"
$ws.Range("E28").Value = "'  -4.77%  "
$ws.Range("D29").Value = "'0.124"
$ws.Range("E29").Value = "'  -6.63%  "
$ws.Range("D30").Value = "'4.128.39"
$ws.Range("E30").Value = "'  -0.05%  "
$ws.Range("D31").Value = "'1.68"
$ws.Range("E31").Value = "'  +2.98%  "
$ws.Range("D32").Value = "'4.10"
$ws.Range("E32").Value = "'  -6.54%  "
$ws.Range("D33").Value = "'0.0564"
$ws.Range("E33").Value = "'  -3.16%  "
$ws.Range("E34").Value = "'  -0.88%  "
$ws.Range("D35").Value = "'4.05"
$ws.Range("E35").Value = "'  -4.05%  "
$ws.Range("D36").Value = "'0.812"
$ws.Range("E36").Value = "'  -11.87%  "
$ws.Range("B37").Value = "'LidoDAOToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.90"
$ws.Range("E37").Value = "'  -5.74%  "
$ws.Range("B38").Value = "'WEMIXToken"
$ws.Range("C38").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.52"
$ws.Range("E38").Value = "'  -23.77%  "
$ws.Range("D39").Value = "'96.46"
$ws.Range("E39").Value = "'  -0.71%  "
$ws.Range("B40").Value = "'Kaspa"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.0660"
$ws.Range("E40").Value = "'  +0.83%  "
$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'16.59"
$ws.Range("E41").Value = "'  -4.86%  "
$ws.Range("D42").Value = "'0.0207"
$ws.Range("E42").Value = "'  -5.23%  "
$ws.Range("D43").Value = "'1.05"
$ws.Range("E43").Value = "'  -5.84%  "
$ws.Range("D44").Value = "'0.0842"
$ws.Range("E44").Value = "'  +13.75%  "
$ws.Range("D45").Value = "'1.270.26"
$ws.Range("E45").Value = "'  -5.28%  "
$ws.Range("D46").Value = "'2.26"
$ws.Range("E46").Value = "'  -7.33%  "
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = "'  -1.30%  "
$ws.Range("E48").Value = "'  -1.83%  "
$ws.Range("D49").Value = "'11.61"
$ws.Range("E49").Value = "'  -1.00%  "
$ws.Range("D50").Value = "'6.20"
$ws.Range("E50").Value = "'  -8.33%  "
$ws.Range("D51").Value = "'41.69"
$ws.Range("E51").Value = "'  -7.79%  "
